$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I, row 2: empty cell, same style as the rest of row 2 (thick bottom border) ---
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

# --- New column I, row 3 (header): 2021, bold/medium-bottom-border like H3 but sz 11 ---
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Font.Size = 11
$ws.Range("I3").Value = 2021

# --- New column I, row 4 (data): 149, same look as H4 but sz 11 ---
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Font.Size = 11
$ws.Range("I4").Value = 149

# --- New column I, row 5 (data, bottom row): 159, same look as H5 but sz 11 ---
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Font.Size = 11
$ws.Range("I5").Value = 159

# Clear clipboard marching ants / match recorded selection in the saved view
$excel.CutCopyMode = $false
$ws.Range("K4").Select()
